$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Coliflor (Macroferia Regional de
# Talca) on 2023-10-16 (Excel serial 45215). It becomes the new row 388,
# pushing every existing record from row 388 down one row (388->389,
# 389->390, ..., 515->516).
#
# Insert a blank row at 388 (this shifts 388..515 down to 389..516),
# duplicate the row that is now at 389 (the old row 388's data) back up
# into 388, then overwrite just the date in the new row with the new
# record's date.

$ws.Rows(388).Insert()

$ws.Range("A389:R389").Copy()
$ws.Range("A388").PasteSpecial()

$ws.Range("D388").Value = 45215
